$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Balogun & Johnson_2004_AMJ_Organizational Restructuring and Middle Manager Sensemaking.pdf", "F1_P11_Balogun & Johnson_2004_AMJ_Organizational Restructuring and Middle Manager Sensemaking.png", 1, 2004, 12, 0, "final_figures/2004/F1_P11_Balogun & Johnson_2004_AMJ_Organizational Restructuring and Middle Manager Sensemaking.png", "process diagram", $true),
    @("Black et al._2004_ASQ_A Dynamic Theory of Expertise_Systems Model.pdf", "F1_P13_Black et al._2004_ASQ_A Dynamic Theory of Expertise_Systems Model.png", 1, 2004, 14, 0, "final_figures/2004/F1_P13_Black et al._2004_ASQ_A Dynamic Theory of Expertise_Systems Model.png", "process diagram", $true),
    @("Ericksen & Dyer_2004_ASQ_Right from the Start- Exploring the Effects of Early Team Events.pdf", "F1_P18_Ericksen & Dyer_2004_ASQ_Right from the Start- Exploring the Effects of Early Team Events.png", 1, 2004, 19, 0, "final_figures/2004/F1_P18_Ericksen & Dyer_2004_ASQ_Right from the Start- Exploring the Effects of Early Team Events.png", "process diagram", $true),
    @("Maitlis & Ozcelik_2004_OrgSci_Toxic Decision Processes.pdf", "F1_P8_Maitlis & Ozcelik_2004_OrgSci_Toxic Decision Processes.png", 1, 2004, 9, 0, "final_figures/2004/F1_P8_Maitlis & Ozcelik_2004_OrgSci_Toxic Decision Processes.png", "process diagram", $true),
    @("Washington_2004_OrgStudies_Field Approaches to Institutional Change.pdf", "F1_P8_Washington_2004_OrgStudies_Field Approaches to Institutional Change.png", 1, 2004, 9, 0, "final_figures/2004/F1_P8_Washington_2004_OrgStudies_Field Approaches to Institutional Change.png", "line graph", $true)
)

$startRow = 28
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}
